# Daily attendance processing - normalize the "Recorded By" (column G) value
# ordering. The attendance sync re-orders the comma-separated list of
# recorder identities into a canonical priority order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Canonical priority order for the "Recorded By" tokens (lower index = sorts first).
# NOTE: comparisons are done with .Equals() to stay case-sensitive, since
# "System" and "system" are distinct tokens here.
$order = @("backup@backdoor.com", "dnasr281@gmail.com", "System", "admin@admin.com", "system")

function Get-Priority($item) {
    for ($i = 0; $i -lt $order.Length; $i++) {
        if ($item.Equals($order[$i])) {
            return $i
        }
    }
    return 999
}

# Stable insertion sort over the priority list above.
function Sort-ByPriority($items) {
    $result = @()
    foreach ($item in $items) {
        $p = Get-Priority $item
        $inserted = $false
        $newResult = @()
        for ($i = 0; $i -lt $result.Length; $i++) {
            $existing = $result[$i]
            $ep = Get-Priority $existing
            $shouldInsert = $false
            if ($inserted -eq $false) {
                if ($p -lt $ep) {
                    $shouldInsert = $true
                }
            }
            if ($shouldInsert -eq $true) {
                $newResult += $item
                $inserted = $true
            }
            $newResult += $existing
        }
        if ($inserted -eq $false) {
            $newResult += $item
        }
        $result = $newResult
    }
    return $result
}

# Find the last used row based on column A (Year).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Text

    if ([string]::IsNullOrEmpty($current)) {
        continue
    }

    $parts = $current.Split(",")
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    $sortedParts = Sort-ByPriority $trimmed
    $newValue = $sortedParts -join ", "

    if (-not $newValue.Equals($current)) {
        $cell.Value = $newValue
    }
}
